$d = $word.ActiveDocument

# Locate the (currently empty) paragraph that immediately precedes the
# "Character will run at constant speed." paragraph. In the source
# document this is an empty <w:p/> right after the "throwGrenade()" line.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "`r" -and $i -lt $d.Paragraphs.Count) {
        $next = $d.Paragraphs.Item($i + 1)
        if ($next.Range.Text -like "Character will run at constant speed.*") {
            $target = $p
            break
        }
    }
}

$newText = "Camera will follow character at a set distance. The camera is position above and behind the character, angled down towards the character"

# Fill the empty paragraph with the new sentence.
$target.Range.Text = $newText

# The document's hidden "_GoBack" bookmark (marking the last edit point)
# needs to move from the final "Enum stance ..." paragraph onto the end
# of this newly typed paragraph, right before its paragraph mark - which
# is where Word leaves it after typing new text.
#
# Adding a bookmark with Bookmarks.Add at a collapsed range sitting
# exactly on "paragraph end - 1" (i.e. immediately before the pilcrow)
# is unreliable, so we temporarily pad the paragraph with two extra
# placeholder characters, anchor the bookmark just before them (a safe,
# non-edge position), then delete the placeholder - leaving the bookmark
# correctly collapsed right after the real text / before the mark.
$pad = "ZZ"
$target.Range.InsertAfter($pad)

$bmPos = $target.Range.End - (1 + $pad.Length)
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$padRange = $d.Range($target.Range.End - (1 + $pad.Length), $target.Range.End - 1)
$padRange.Text = ""
